$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: clear the stray/garbled row (values only, keep formatting) ---
$ws.Range("A3:F3").ClearContents() | Out-Null
$ws.Rows.Item(3).RowHeight = 15.75

# --- Row 8: clear the duplicate row (values only, keep formatting) ---
$ws.Range("A8:F8").ClearContents() | Out-Null
$ws.Rows.Item(8).RowHeight = 15.75

# --- Row 10: fill in a new recharge record (Annual Lite plan) ---
# Match the date column's formatting (copy from D9) before writing the date value.
$ws.Range("D9").Copy() | Out-Null
$ws.Range("D10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A10").Value = 108
$ws.Range("B10").Value = 399
$ws.Range("C10").Value = "postpaid"
$ws.Range("D10").Value = 45665
$ws.Range("E10").Value = "success"
$ws.Range("F10").Value = "Annual Lite"
$ws.Rows.Item(10).RowHeight = 29.25

# --- Row 12: fill in a new recharge record (Monthly Plan) ---
$ws.Range("A12").Value = 110
$ws.Range("B12").Value = 299
$ws.Range("C12").Value = "postpaid"
$ws.Range("D12").Value = 45667
$ws.Range("E12").Value = "success"
$ws.Range("F12").Value = "Monthly Plan"
$ws.Rows.Item(12).RowHeight = 29.25

# --- Row 15: now a two-line wrapped row like its neighbours ---
$ws.Rows.Item(15).RowHeight = 29.25

# --- Column widths (now that plan_name/date text is wider) ---
$ws.Columns.Item(1).ColumnWidth = 10.666666666666666
$ws.Columns.Item(2).ColumnWidth = 11.0
$ws.Columns.Item(3).ColumnWidth = 16.666666666666668
$ws.Columns.Item(4).ColumnWidth = 19.333333333333332
$ws.Columns.Item(5).ColumnWidth = 11.833333333333334

# --- Selection left where the author last clicked ---
$ws.Range("G10").Select() | Out-Null
